$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper "template" ranges, taken from the ORIGINAL (still untouched at this
# point) paragraph 1 ("Hi, this is my first {@define jamal="), which already
# carries correctly-serialised run formatting (<w:lang w:val="en-US"/>, and
# for "Hi" also <w:i/><w:iCs/>). Copying their FormattedText and then
# re-setting .Text on the pasted copy reliably reproduces that run
# formatting on brand-new text, which plain Font.Italic / LanguageID
# assignment does not always do for runs created in a previously-empty
# paragraph.
# ---------------------------------------------------------------------------
$italicTemplate = $d.Range(0, 2)    # "Hi"  -> i / iCs / lang=en-US
$plainTemplate  = $d.Range(35, 36)  # "="   -> lang=en-US

$italicFT = $italicTemplate.FormattedText
$plainFT  = $plainTemplate.FormattedText

# ---------------------------------------------------------------------------
# Paragraph 1: "Hi, this is my first {@define jamal="
#   -> "{@ident {jama mamal}}Hi, this is my first {@define jamal=0"
# ---------------------------------------------------------------------------

$p1Start = $d.Paragraphs.Item(1).Range.Start
$insHead = $d.Range($p1Start, $p1Start)
$insHead.FormattedText = $italicFT
$headRange = $d.Range($p1Start, $p1Start + 2)
$headRange.Text = "{@ident {jama mamal}}"

$p1End = $d.Paragraphs.Item(1).Range.End - 1
$insTail = $d.Range($p1End, $p1End)
$insTail.FormattedText = $plainFT
$tailRange = $d.Range($p1End, $p1End + 1)
$tailRange.Text = "0"

# ---------------------------------------------------------------------------
# Paragraph 2: "" -> "1"
# ---------------------------------------------------------------------------
$p2Start = $d.Paragraphs.Item(2).Range.Start
$ins2 = $d.Range($p2Start, $p2Start)
$ins2.FormattedText = $plainFT
$r2 = $d.Range($p2Start, $p2Start + 1)
$r2.Text = "1"

# ---------------------------------------------------------------------------
# Paragraph 3: "" -> "2"
# ---------------------------------------------------------------------------
$p3Start = $d.Paragraphs.Item(3).Range.Start
$ins3 = $d.Range($p3Start, $p3Start)
$ins3.FormattedText = $plainFT
$r3 = $d.Range($p3Start, $p3Start + 1)
$r3.Text = "2"

# Paragraph 4 ("Jamal") is unchanged.

# ---------------------------------------------------------------------------
# Paragraph 5: "" -> "3"
# ---------------------------------------------------------------------------
$p5Start = $d.Paragraphs.Item(5).Range.Start
$ins5 = $d.Range($p5Start, $p5Start)
$ins5.FormattedText = $plainFT
$r5 = $d.Range($p5Start, $p5Start + 1)
$r5.Text = "3"

# Paragraph 6 ("}{jamal} processed document.") is unchanged.
# Paragraph 7 ("") is unchanged.

# ---------------------------------------------------------------------------
# Paragraph 8: "Wouw" -> "Wouw{jamal}"
#   "{" is italic, "jamal}" is plain.
# ---------------------------------------------------------------------------
$p8End = $d.Paragraphs.Item(8).Range.End - 1

$insBrace = $d.Range($p8End, $p8End)
$insBrace.FormattedText = $italicFT
$braceRange = $d.Range($p8End, $p8End + 2)
$braceRange.Text = "{"

$p8End2 = $d.Paragraphs.Item(8).Range.End - 1
$insRest = $d.Range($p8End2, $p8End2)
$insRest.FormattedText = $plainFT
$restRange = $d.Range($p8End2, $p8End2 + 1)
$restRange.Text = "jamal}"

Write-Output "done"
